$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "34.627.23"
$ws.Range("E2").Value = "  +1.28%  "

# Row 3
$ws.Range("D3").Value = "1.796.37"
$ws.Range("E3").Value = "  +0.80%  "

# Row 4
$ws.Range("E4").Value = "  -0.11%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.10"
$ws.Range("E5").Value = "  +0.55%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.559"
$ws.Range("E6").Value = "  +2.29%  "

# Row 7
$ws.Range("E7").Value = "  -0.12%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "33.02"
$ws.Range("E8").Value = "  +4.14%  "

# Row 9
$ws.Range("E9").Value = "  +2.05%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0696"
$ws.Range("E10").Value = "  +1.33%  "

# Row 11
$ws.Range("E11").Value = "  +0.54%  "

# Row 12
$ws.Range("D12").Value = "2.054.37"
$ws.Range("E12").Value = "  +0.81%  "

# Row 13
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.10"
$ws.Range("E13").Value = "  +0.99%  "

# Row 14
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.796.24"
$ws.Range("E14").Value = "  +0.96%  "

# Row 15
$ws.Range("E15").Value = "  +2.48%  "

# Row 16
$ws.Range("D16").Value = "34.558.31"
$ws.Range("E16").Value = "  +1.26%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.30"
$ws.Range("E17").Value = "  +3.06%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.93"
$ws.Range("E18").Value = "  +1.55%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "248.10"
$ws.Range("E19").Value = "  +0.83%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0802"
$ws.Range("E20").Value = "  +1.56%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.30"
$ws.Range("E21").Value = "  +2.98%  "

# Row 22
$ws.Range("E22").Value = "  -0.14%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.19"
$ws.Range("E23").Value = "  +2.25%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.08"
$ws.Range("E24").Value = "  +1.52%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "166.43"
$ws.Range("E25").Value = "  +2.76%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.31"
$ws.Range("E26").Value = "  +1.98%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.59"
$ws.Range("E27").Value = "  +1.71%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.117"
$ws.Range("E28").Value = "  +2.50%  "

# Row 29
$ws.Range("E29").Value = "  -0.23%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.11"
$ws.Range("E30").Value = "  +12.86%  "

# Row 31
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0526"
$ws.Range("E31").Value = "  +1.55%  "

# Row 32
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.24"
$ws.Range("E32").Value = "  +1.05%  "

# Row 33
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.82"
$ws.Range("E33").Value = "  +2.63%  "

# Row 34
$ws.Range("E34").Value = "  +2.69%  "

# Row 35
$ws.Range("D35").Value = "1.429.16"
$ws.Range("E35").Value = "  -0.99%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.59"
$ws.Range("E36").Value = "  +7.79%  "

# Row 37
$ws.Range("E37").Value = "  +3.15%  "

# Row 38
$ws.Range("E38").Value = "  +1.11%  "

# Row 39
$ws.Range("E39").Value = "  +2.05%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "85.73"
$ws.Range("E40").Value = "  +6.76%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.42"
$ws.Range("E41").Value = "  +1.76%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.935"
$ws.Range("E42").Value = "  +1.30%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.76"
$ws.Range("E43").Value = "  +3.14%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.66"
$ws.Range("E44").Value = "  +0.49%  "

# Row 45
$ws.Range("E45").Value = "  +3.63%  "

# Row 46
$ws.Range("E46").Value = "  +1.16%  "

# Row 47
$ws.Range("E47").Value = "  +0.39%  "

# Row 48
$ws.Range("D48").Value = "1.953.92"
$ws.Range("E48").Value = "  +0.70%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "106.27"
$ws.Range("E49").Value = "  +1.51%  "

# Row 50
$ws.Range("B50").Value = "PaxDollar"
$ws.Range("C50").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  -0.11%  "

# Row 51
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₆0129"
$ws.Range("E51").Value = "  -5.97%  "
